$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 230, shifting existing rows 230+ down by one.
$ws.Rows("230:230").Insert()

# Populate the newly inserted row 230 with the new data.
$ws.Cells.Item(230, 1).Value = 10
$ws.Cells.Item(230, 2).Value = "Vega Modelo de Temuco"
$ws.Cells.Item(230, 3).Value = "La Araucanía"
$ws.Cells.Item(230, 4).Value = 44510
$ws.Cells.Item(230, 5).Value = 9
$ws.Cells.Item(230, 6).Value = 100112043
$ws.Cells.Item(230, 7).Value = "Pepino ensalada"
$ws.Cells.Item(230, 8).Value = "Sin especificar"
$ws.Cells.Item(230, 9).Value = "Primera"
$ws.Cells.Item(230, 10).Value = 100
$ws.Cells.Item(230, 11).Value = 10000
$ws.Cells.Item(230, 12).Value = 10000
$ws.Cells.Item(230, 13).Value = 10000
$ws.Cells.Item(230, 14).Value = "$/caja 60 unidades"
$ws.Cells.Item(230, 15).Value = "Región de O'Higgins"
$ws.Cells.Item(230, 16).Value = 167
$ws.Cells.Item(230, 17).Value = 60
$ws.Cells.Item(230, 18).Value = "Hortaliza"

# Ensure the date cell keeps the same date/time number format as the rest of column D.
$ws.Cells.Item(230, 4).NumberFormat = $ws.Cells.Item(229, 4).NumberFormat
